$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 40 (Segunda, 2021-12-02 week data, 2600/2600/2600) which is
# superseded by the new weekly rows being added at the top of this date block.
$ws.Rows.Item(40).Delete()

# Insert two new rows at the top of the block (rows 37-38) for the new week
# (2021-10-20 -> serial 44489), pushing the remaining rows down by one.
$ws.Range("37:38").Insert()

# Row 37: Espárragos, Sin especificar, Banquete
$ws.Cells.Item(37, 1).Value = 12
$ws.Cells.Item(37, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(37, 3).Value = "Metropolitana"
$ws.Cells.Item(37, 4).Value = 44489
$ws.Cells.Item(37, 5).Value = 13
$ws.Cells.Item(37, 6).Value = 300000000
$ws.Cells.Item(37, 7).Value = "Espárragos"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Banquete"
$ws.Cells.Item(37, 10).Value = 330
$ws.Cells.Item(37, 11).Value = 1200
$ws.Cells.Item(37, 12).Value = 1200
$ws.Cells.Item(37, 13).Value = 1200
$ws.Cells.Item(37, 14).Value = "$/kilo"
$ws.Cells.Item(37, 15).Value = "Región Metropolitana"
$ws.Cells.Item(37, 16).Value = 1200
$ws.Cells.Item(37, 17).Value = 1
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Row 38: Espárragos, Sin especificar, Primera
$ws.Cells.Item(38, 1).Value = 12
$ws.Cells.Item(38, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(38, 3).Value = "Metropolitana"
$ws.Cells.Item(38, 4).Value = 44489
$ws.Cells.Item(38, 5).Value = 13
$ws.Cells.Item(38, 6).Value = 300000000
$ws.Cells.Item(38, 7).Value = "Espárragos"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 450
$ws.Cells.Item(38, 11).Value = 1000
$ws.Cells.Item(38, 12).Value = 1000
$ws.Cells.Item(38, 13).Value = 1000
$ws.Cells.Item(38, 14).Value = "$/kilo"
$ws.Cells.Item(38, 15).Value = "Región Metropolitana"
$ws.Cells.Item(38, 16).Value = 1000
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"
